# Fix partner list report: update header labels and refresh the saved
# worksheet selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text fixes
$ws.Range("A2").Value = "Is Customer"
$ws.Range("A3").Value = "Is Supplier"
$ws.Range("B5").Value = "Partner Name"
$ws.Range("AB5").Value = "Bank Name"

# Refresh the saved view: active cell AB6, top-left visible cell Y1
$ws.Range("AB6").Select()
$excel.ActiveWindow.ScrollColumn = 25
$excel.ActiveWindow.ScrollRow = 1
